# EDM-5: crud functionality implemented, patch added
# EDM-5: UI is divided into components, improving business logic
#
# Rewrites the "Документы" (Documents) list: the old row 1 record is removed,
# the old row 2 record is replaced, and two more records are written into
# rows 3-4 (re-using the previously-empty placeholder rows). The trailing
# placeholder rows keep their date number-format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal text even when it looks like
# an ambiguous date (e.g. day <= 12), so Excel does not silently convert it
# into a date serial number. We briefly mark the cell as Text, assign the
# value, then restore the real (date) display format.
function Set-TextCellValue($range, $value, $numberFormat) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $numberFormat
}

# ---------------------------------------------------------------------
# 1) Row 1 (id=1, "Доверенность" ...) is deleted entirely.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).Clear()

# ---------------------------------------------------------------------
# 2) Row 2 becomes: 1 | 1232 | gfgfgf | 13.01.2024 | 20.01.2024 | 2
#    (D2/E2 already carry a dd.mm.yyyy style and neither value is an
#    ambiguous day, so a plain value assignment keeps them as text.)
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1232
$ws.Range("C2").Value = "gfgfgf"
$ws.Range("D2").Value = "13.01.2024"
$ws.Range("E2").Value = "20.01.2024"
$ws.Range("F2").Value = 2

# ---------------------------------------------------------------------
# 3) Row 3 becomes: 2 | 13123 | Абоба | 12.01.2024 | 04.02.2024 | 1
#    Both dates have a day-of-month <= 12, which Excel would otherwise
#    read as an (ambiguous) date and silently convert to a serial
#    number, so they are written through the text-safe helper.
# ---------------------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 13123
$ws.Range("C3").Value = "Абоба"
Set-TextCellValue $ws.Range("D3") "12.01.2024" "dd.mm.yyyy"
Set-TextCellValue $ws.Range("E3") "04.02.2024" "dd.mm.yyyy"
$ws.Range("F3").Value = 1

# ---------------------------------------------------------------------
# 4) Row 4 becomes: 3 | 1 | fgfgfgfgg | 14.01.2021 | 28.01.2024 | 6
# ---------------------------------------------------------------------
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "fgfgfgfgg"
$ws.Range("D4").Value = "14.01.2021"
$ws.Range("D4").NumberFormat = "dd.mm.yyyy"
$ws.Range("E4").Value = "28.01.2024"
$ws.Range("E4").NumberFormat = "dd.mm.yyyy"
$ws.Range("F4").Value = 6

# ---------------------------------------------------------------------
# 5) Rows 5-8 stay as empty, date-formatted placeholder rows (unchanged).
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 6) Update the active selection shown in the workbook.
# ---------------------------------------------------------------------
[void]$ws.Range("G6").Select()
